$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.420949295932613
$ws.Range("B2").Value = -4.257060648907316

$ws.Range("A3").Value = -0.5046424108579408
$ws.Range("B3").Value = 1.188682204638741

$ws.Range("A4").Value = 1.006628774259225
$ws.Range("B4").Value = -3.096770889681989

$ws.Range("A5").Value = 0.6790302894739064
$ws.Range("B5").Value = 0.4022692427991631

$ws.Range("A6").Value = -0.8155874771776408
$ws.Range("B6").Value = -1.881942139627402

$ws.Range("A7").Value = -0.0612483592114824
$ws.Range("B7").Value = -0.697440638001632

$ws.Range("A8").Value = 0.80421956980311
$ws.Range("B8").Value = 0.7631689658383464

$ws.Range("A9").Value = 0.3114309331481631
$ws.Range("B9").Value = 1.02485014504204

$ws.Range("A10").Value = -0.1775653614068756
$ws.Range("B10").Value = -2.526891038649877

$ws.Range("A11").Value = 0.3461737102865866
$ws.Range("B11").Value = -0.587349174488855

$ws.Range("A12").Value = 0.5244127092989271
$ws.Range("B12").Value = 0.6304501353211588

$ws.Range("A13").Value = 0.4048370064451582
$ws.Range("B13").Value = 0.06255283358463837
